# New data for signup is added
$wb = $excel.ActiveWorkbook

# 1. Add a new "Signup" sheet and move it to the front of the workbook.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Signup"
$newSheet.Move($wb.Worksheets.Item(1))
# NOTE: .Move() repositions sheets by index under the hood, so any handle
# obtained before the move can resolve to the wrong sheet afterwards.
# Re-fetch a fresh handle by name for all subsequent writes.
$signup = $wb.Worksheets.Item("Signup")

# 2. Populate the Signup sheet with the same layout style as the other
#    sheets: a numeric header row, a section-title row, a field-name row
#    and a data row.
$signup.Range("A1").Value = 0
$signup.Range("B1").Value = 1
$signup.Range("C1").Value = 2
$signup.Range("D1").Value = 3
$signup.Range("E1").Value = 4
$signup.Range("F1").Value = 5
$signup.Range("G1").Value = 6
$signup.Range("H1").Value = 7
$signup.Range("I1").Value = 8

$signup.Range("A2").Value = "Signup"

$signup.Range("A3").Value = "RunMode"
$signup.Range("B3").Value = "Firstname"
$signup.Range("C3").Value = "Lastname"
$signup.Range("D3").Value = "Address"
$signup.Range("E3").Value = "Phone"
$signup.Range("F3").Value = "email"
$signup.Range("G3").Value = "Password"
$signup.Range("H3").Value = "City"
$signup.Range("I3").Value = "Zipcode"

$signup.Range("A4").Value = "Signup"
$signup.Range("B4").Value = "Mason"
$signup.Range("C4").Value = "williams"
$signup.Range("D4").Value = "6136 Walraven Cir"
$signup.Range("E4").Value = "966-659-7666"
$signup.Range("F4").Value = "jamesrobert441@gmail.com"
$signup.Range("G4").Value = 123456
$signup.Range("H4").Value = "Memphis"
$signup.Range("I4").Value = 38119

# 3. Column widths for the new sheet.
$signup.Columns.Item(2).ColumnWidth = 20.7109375
$signup.Columns.Item(3).ColumnWidth = 10.85546875
$signup.Columns.Item(4).ColumnWidth = 18.85546875
$signup.Columns.Item(5).ColumnWidth = 12.85546875
$signup.Columns.Item(6).ColumnWidth = 27.5703125

# 4. The existing Login sheet: the FbPasswod demo value now reuses the
#    mani6747@gmail.com sample address instead of ebiztesting8@gmail.com.
$login = $wb.Worksheets.Item("Login")
$login.Range("B4").Value = "mani6747@gmail.com"

Write-Output "done"
